$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 0.5778618378969409
$ws.Range("D2").Value = 0.1576473297629773
$ws.Range("E2").Value = "norm_coldread_gaze_wpm_median"
$ws.Range("F2").Value = "nan"

# Row 3
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 0.6087350244037076
$ws.Range("D3").Value = 0.1210832617488018
$ws.Range("E3").Value = "norm_coldread_gaze_wpm_median"
$ws.Range("F3").Value = "nan"

# Row 4
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 0.5615040311219779
$ws.Range("D4").Value = 0.09119998029116146
$ws.Range("E4").Value = "norm_coldread_gaze_wpm_median"
$ws.Range("F4").Value = "nan"
